$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix precision of existing A53 timestamp value
$ws.Range("A53").Value = 44366.76719248958

# Add new row 54 data
$ws.Range("A54").Value = 44367.76684843018
$ws.Range("B54").Value = 77874
$ws.Range("C54").Value = 65506
$ws.Range("D54").Value = 3480
$ws.Range("E54").Value = 2098
$ws.Range("F54").Value = 1484
$ws.Range("G54").Value = 20631
$ws.Range("H54").Value = 1453
$ws.Range("I54").Value = 893
$ws.Range("J54").Value = 184

# Match the date-time style used by column A (style index 2)
$ws.Range("A54").NumberFormat = $ws.Range("A53").NumberFormat
